$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank memory cells for Tree BFS/DFS 10000 (columns G/H) for rows 19-21
# (these cells previously only carried a style with no value; clear that
# formatting back to the default style as part of giving them real values)
$ws.Range("G19:H21").ClearFormats()

$ws.Range("G19").Value = 0.1154601
$ws.Range("H19").Value = 0.1431686

$ws.Range("G20").Value = 0.1081391
$ws.Range("H20").Value = 0.1331268

$ws.Range("G21").Value = 0.0997361
$ws.Range("H21").Value = 0.1361017

# Add missing memory values for row 26 (Tree BFS/DFS 1000, columns I/J)
$ws.Range("I26").Value = 0.0231459
$ws.Range("J26").Value = 0.0147596

# Add a new row 28 of data (Tree BFS/DFS 50000, columns K/L)
$ws.Range("K28").Value = 0.6117588
$ws.Range("L28").Value = 0.5543064
